# Applies the cryptos.xlsx price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.425.56"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'334.37"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "'0.4551"
$ws.Range("E8").Value = "  +5.53%  "
$ws.Range("D9").Value = "'53.16"
$ws.Range("E9").Value = "  +15.29%  "
$ws.Range("D10").Value = "'0.08932"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").Value = "'2.106.11"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'6.847"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "'8.067"
$ws.Range("E15").Value = "  +4.98%  "
$ws.Range("D16").Value = "'96.58"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "'0.00001141"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'0.06651"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'6.338"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'30.488.22"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'2.364"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").Value = "'2.353.97"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'22.28"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'2.536"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "'162.70"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").Value = "'133.11"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'1.212"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").Value = "'1.667"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "'6.386"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").Value = "'3.940"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").Value = "'10.46"
$ws.Range("E36").Value = "  +5.92%  "
$ws.Range("D37").Value = "'5.791"
$ws.Range("E37").Value = "  +6.26%  "
$ws.Range("D38").Value = "'0.02590"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "'0.06844"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("D40").Value = "'0.2301"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").Value = "'12.71"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'0.6889"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "'1.250"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'2.319"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.06"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").Value = "'0.6377"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +21.25%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'83.45"
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "'0.3391"
$ws.Range("E51").Value = "  +23.30%  "
